$d = $word.ActiveDocument

# The document currently ends with:
#   ...VL7<bookmark _GoBack>
#   -----------------------------------------------------------------------------------------------------------------------------------
#   <empty paragraph>
#
# It needs to end with:
#   ...VL7
#   -----------------------------------------------------------------------------------------------------------------------------------
#   Praktikumsblatt 6(aber etwas abgewandelter als wir es gelöst haben)<bookmark _GoBack>
#
# i.e. the trailing empty paragraph gets the new sentence, and the
# "_GoBack" bookmark moves from the "VL7" paragraph to the end of that
# new (now non-empty) last paragraph.

$newText = "Praktikumsblatt 6(aber etwas abgewandelter als wir es gelöst haben)"

# The trailing empty paragraph is the very last paragraph of the body.
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)

# Put a collapsed "_GoBack" bookmark at the (currently empty) paragraph's
# position first. Adding a bookmark named "_GoBack" automatically removes
# it from wherever it previously lived (the "VL7" paragraph), matching
# Word's own single-instance "_GoBack" behaviour. Because the bookmark is
# added before the new text is typed, the inserted text lands in front of
# the (still collapsed) bookmark - exactly like a user typing at that spot,
# and leaving the bookmark collapsed at the end of the paragraph.
$d.Bookmarks.Add("_GoBack", $lastPara.Range)

# Now fill in the paragraph's text.
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.Text = $newText
